# Append two new "LeetCode problem" paragraphs right after the
# "- Solution is dynamic programming" paragraph, at the end of the document:
#
#   - Merge Intervals (Didn't remember this problem, was stumped...)
#       (two runs: "- Merge Intervals" and " (Didn't remember this
#        problem, was stumped...)" -- the latter keeps its leading space,
#        hence xml:space="preserve")
#   - Reconstruct Itinerary
#
$d = $word.ActiveDocument

# Locate the anchor paragraph ("- Solution is dynamic programming") by
# scanning Paragraphs (more robust here than narrowing via Find, whose
# resulting Range.Paragraphs(1) does not reliably expand back to the full
# containing paragraph).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "- Solution is dynamic programming") {
        $anchorIndex = $i
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph '- Solution is dynamic programming'"
}
$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Create a fresh empty paragraph right after the anchor; we'll fill it (and
# the one after it) via raw WordOpenXML so the " (Didn't remember ...)" text
# lands in its own run exactly as in the source edit.
$null = $anchorPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 1)

$apostrophe = [string]([char]0x2019)
$ellipsis = [string]([char]0x2026)

$mergeIntervalsRun2 = " (Didn" + $apostrophe + "t remember this problem, was stumped" + $ellipsis + ")"

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:r><w:t>- Merge Intervals</w:t></w:r><w:r><w:t xml:space="preserve">' + $mergeIntervalsRun2 + '</w:t></w:r></w:p>' +
  '<w:p><w:r><w:t>- Reconstruct Itinerary</w:t></w:r></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$null = $newPara.Range.InsertXML($xml)

# InsertXML inserted the two brand-new paragraphs ahead of $newPara's own
# paragraph mark, leaving that original (now-redundant) empty paragraph
# dangling at the very end -- remove it (and its mark) so the document ends
# with "...- Reconstruct Itinerary" immediately followed by sectPr, same as
# the source edit.
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
if ($trailing.Range.Text.Trim() -eq "") {
    $cleanupRange = $d.Range($trailing.Range.Start - 1, $trailing.Range.End)
    $null = $cleanupRange.Delete()
}

$check1 = $d.Paragraphs.Item($anchorIndex + 1).Range.Text
$check2 = $d.Paragraphs.Item($anchorIndex + 2).Range.Text
if (-not ($check1 -match "^- Merge Intervals \(Didn.t remember this problem, was stumped.\)")) {
    throw "Unexpected text for 'Merge Intervals' paragraph: [$check1]"
}
if (-not ($check2 -match "^- Reconstruct Itinerary")) {
    throw "Unexpected text for 'Reconstruct Itinerary' paragraph: [$check2]"
}

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
